$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns, matching the style used by the other
# header cells (row 1, e.g. H1). Copy H1 (value + formatting) into I1/J1
# first so the new cells reuse the existing bold/centered/bordered style,
# then overwrite the copied text with the real header labels.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-21: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
